$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1049
    $ws.Range("F3").Value = 430
    $ws.Range("F4").Value = 3167
    $ws.Range("F6").Value = 642
}
